$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 626  # was 621
$ws.Range("F7").Value = 468  # was 464
$ws.Range("F8").Value = 252  # was 248
$ws.Range("F9").Value = 1127  # was 1119
$ws.Range("F11").Value = 171  # was 167
$ws.Range("F12").Value = 60  # was 58
$ws.Range("F13").Value = 754  # was 753
$ws.Range("F14").Value = 395  # was 391
$ws.Range("F15").Value = 54  # was 53
$ws.Range("F17").Value = 195  # was 193
$ws.Range("F18").Value = 8  # was 7
$ws.Range("F19").Value = 382  # was 380
$ws.Range("F20").Value = 6192  # was 6179
$ws.Range("F22").Value = 54  # was 53
$ws.Range("F24").Value = 7172  # was 7159
$ws.Range("F27").Value = 3279  # was 3277
$ws.Range("F28").Value = 402  # was 396
$ws.Range("F29").Value = 786  # was 785
$ws.Range("F30").Value = 4470  # was 4471
$ws.Range("F32").Value = 151  # was 149
$ws.Range("F33").Value = 155  # was 154
$ws.Range("F34").Value = 1233  # was 1223
$ws.Range("F38").Value = 977  # was 970
$ws.Range("F39").Value = 1251  # was 1241
$ws.Range("F40").Value = 2070  # was 2066

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 233  # was 232
$ws.Range("F3").Value = 1166  # was 1165

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 233  # was 232
$ws.Range("F4").Value = 1166  # was 1165
$ws.Range("F9").Value = 626  # was 621
$ws.Range("F10").Value = 468  # was 464
$ws.Range("F11").Value = 252  # was 248
$ws.Range("F12").Value = 1127  # was 1119
$ws.Range("F14").Value = 171  # was 167
$ws.Range("F15").Value = 60  # was 58
$ws.Range("F16").Value = 754  # was 753
$ws.Range("F17").Value = 395  # was 391
$ws.Range("F18").Value = 54  # was 53
$ws.Range("F21").Value = 195  # was 193
$ws.Range("F22").Value = 8  # was 7
$ws.Range("F23").Value = 382  # was 380
$ws.Range("F24").Value = 6192  # was 6179
$ws.Range("F25").Value = 6192  # was 6179
$ws.Range("F27").Value = 54  # was 53
$ws.Range("F29").Value = 7172  # was 7159
$ws.Range("F32").Value = 3279  # was 3277
$ws.Range("F33").Value = 402  # was 396
$ws.Range("F34").Value = 786  # was 785
$ws.Range("F35").Value = 4470  # was 4471
$ws.Range("F38").Value = 151  # was 149
$ws.Range("F39").Value = 155  # was 154
$ws.Range("F40").Value = 1233  # was 1223
$ws.Range("F44").Value = 977  # was 970
$ws.Range("F45").Value = 1251  # was 1241
$ws.Range("F47").Value = 2070  # was 2066
